# Apply "break out stock.yaml completed" edit to SBICARD.NS 1wk sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a couple of previously-wrong values in existing rows ---
$ws.Cells.Item(58, 17).Value = 0     # Q58: detect_structure 2 -> 0
$ws.Cells.Item(222, 15).Value = 2    # O222: isPivot 0 -> 2

# --- Backfill the trailing "backup" column for the last two existing rows ---
$ws.Cells.Item(224, 18).Value = 0    # R224: inlineStr "" -> numeric 0
$ws.Cells.Item(225, 18).Value = 0    # R225: inlineStr "" -> numeric 0

# --- Append seven new weekly rows (226-232) ---
$newRows = @(
    @(45474, 724,               729.7000122070312, 708.75,            721.9500122070312, 721.9500122070312, 6043031,  2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
    @(45481, 723.25,            750,               714.3499755859375, 738.6500244140625, 738.6500244140625, 8051824,  2024, 7, 8,  0, 0, 0, 28, 1, 0, 0),
    @(45488, 740.4500122070312, 743.75,            716.2999877929688, 718.5999755859375, 718.5999755859375, 3921653,  2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 717,               745,               711.5999755859375, 721.7000122070312, 721.7000122070312, 10286468, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 703,               727.75,            702.25,            714.5499877929688, 714.5499877929688, 9724322,  2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 708,               721.9500122070312, 697.4500122070312, 709.7999877929688, 709.7999877929688, 4166685,  2024, 8, 5,  0, 0, 0, 32, 0, 0, 0),
    @(45516, 710,               710.9500122070312, 689,               698.6500244140625, 698.6500244140625, 3449061,  2024, 8, 12, 0, 0, 0, 33, 0, 0, 0)
)

$r = 226
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    # R column (backup) is left blank/empty, matching the diff's inlineStr placeholder.
    $r = $r + 1
}
